# Add a new "points" column (E) to Sheet1, with a header cell styled like
# the existing header row (bold, centered, top-aligned) but with a
# left+right thin border instead of a full box, and fill the column with
# a repeating 1-2-3-4-5 pattern for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell E1 -------------------------------------------------
$header = $ws.Cells.Item(1, 5)
$header.Value = "points"
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.Item(7).Weight = 2    # xlEdgeLeft, xlThin
$header.Borders.Item(10).Weight = 2   # xlEdgeRight, xlThin

# --- Data rows 2..104: repeating 1,2,3,4,5 pattern -------------------
$points = 1
for ($row = 2; $row -le 104; $row++) {
    $ws.Cells.Item($row, 5).Value = $points
    $points = $points + 1
    if ($points -gt 5) {
        $points = 1
    }
}

# --- Selection matches the post-edit state in the author's session ---
$null = $ws.Range("F104").Select()

Write-Host "points column added"
